$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above the current row 445 ("1a (guarda)" / 44421 entry)
# to make room for a new weekly data pair (44509), shifting all subsequent
# rows down by two (old 445 -> 447, ..., old 462 -> 464).
$ws.Rows("445:446").Insert()

# New row 445: "1a (cosecha)" entry for 44509 (Region de Arica y Parinacota)
$ws.Cells.Item(445, 1).Value = 8
$ws.Cells.Item(445, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(445, 3).Value = "Coquimbo"
$ws.Cells.Item(445, 4).Value = 44509
$ws.Cells.Item(445, 5).Value = 4
$ws.Cells.Item(445, 6).Value = 100112004
$ws.Cells.Item(445, 7).Value = "Cebolla"
$ws.Cells.Item(445, 8).Value = "Sin especificar"
$ws.Cells.Item(445, 9).Value = "1a (cosecha)"
$ws.Cells.Item(445, 10).Value = 3000
$ws.Cells.Item(445, 11).Value = 4800
$ws.Cells.Item(445, 12).Value = 5000
$ws.Cells.Item(445, 13).Value = 4900
$ws.Cells.Item(445, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(445, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(445, 16).Value = 272
$ws.Cells.Item(445, 17).Value = 18
$ws.Cells.Item(445, 18).Value = "Hortaliza"

# New row 446: "2a (cosecha)" entry for 44509 (Region de Arica y Parinacota)
$ws.Cells.Item(446, 1).Value = 8
$ws.Cells.Item(446, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(446, 3).Value = "Coquimbo"
$ws.Cells.Item(446, 4).Value = 44509
$ws.Cells.Item(446, 5).Value = 4
$ws.Cells.Item(446, 6).Value = 100112004
$ws.Cells.Item(446, 7).Value = "Cebolla"
$ws.Cells.Item(446, 8).Value = "Sin especificar"
$ws.Cells.Item(446, 9).Value = "2a (cosecha)"
$ws.Cells.Item(446, 10).Value = 1640
$ws.Cells.Item(446, 11).Value = 4500
$ws.Cells.Item(446, 12).Value = 4600
$ws.Cells.Item(446, 13).Value = 4550
$ws.Cells.Item(446, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(446, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(446, 16).Value = 253
$ws.Cells.Item(446, 17).Value = 18
$ws.Cells.Item(446, 18).Value = "Hortaliza"
